# Week 17 data logging for 49ers 2021 Team Data
# Appends this week's per-play yardage samples to the running logs on the
# YDS sheet, updates the season cumulative totals on OFF / DEF / ST /
# TURNS / PEN, and extends the special-teams kick-distance logs on ST.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append this week's individual play-yardage samples to the
# running space-delimited logs.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 2 -2 6 6 2 13 5 1 1 13 16 6 7 -1 0 2 2 6 3 6 1 10 3 6 6 1 8 2 -2 0 37 3 4 5 2 2"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 3 3 12 12 12 1 6 17 27 12 43 6 8 13 45 29"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " -4 6 3 -1 2 6 4 2 5 3 6 5 -1 2 6 4 2 4 3 2 8 2 8 -1 5 3 2"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " 6 3 9 -2 2 8 6 15 24 0 8 6 9 4 8 3 20 9 5 6 14"

# ---------------------------------------------------------------------
# OFF sheet - season cumulative offensive totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 227
$offWs.Range("E2").Value = 10
$offWs.Range("F2").Value = 67
$offWs.Range("G2").Value = 62
$offWs.Range("I2").Value = 10
$offWs.Range("J2").Value = 34
$offWs.Range("L2").Value = 240
$offWs.Range("M2").Value = 157
$offWs.Range("O2").Value = 20
$offWs.Range("Q2").Value = 496

$offWs.Range("C3").Value = 156
$offWs.Range("E3").Value = 33
$offWs.Range("F3").Value = 93
$offWs.Range("H3").Value = 31
$offWs.Range("I3").Value = 53
$offWs.Range("J3").Value = 37
$offWs.Range("N3").Value = 17

# ---------------------------------------------------------------------
# DEF sheet - season cumulative defensive totals (Home row 2, Road row 3)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 186
$defWs.Range("F2").Value = 64
$defWs.Range("G2").Value = 54
$defWs.Range("H2").Value = 6
$defWs.Range("J2").Value = 27
$defWs.Range("L2").Value = 249
$defWs.Range("M2").Value = 167
$defWs.Range("O2").Value = 24
$defWs.Range("Q2").Value = 472

$defWs.Range("B3").Value = 11
$defWs.Range("C3").Value = 175
$defWs.Range("E3").Value = 36
$defWs.Range("F3").Value = 93
$defWs.Range("G3").Value = 33
$defWs.Range("H3").Value = 30
$defWs.Range("I3").Value = 60
$defWs.Range("J3").Value = 61
$defWs.Range("N3").Value = 16

# ---------------------------------------------------------------------
# ST sheet - season cumulative special-teams totals plus the per-kick
# distance logs (row4 = distance, row5 = return, row6 = return man)
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 82
$stWs.Range("D2").Value = 51
$stWs.Range("F2").Value = 584
$stWs.Range("G2").Value = 569
$stWs.Range("J2").Value = 289
$stWs.Range("K2").Value = 274

$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 52 28 36 43"
$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 66 65 61 66 57"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 15 0 8 0"
$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 21 18 26 20 13"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0 9 0 18 11 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Text + " 28 7"

# ---------------------------------------------------------------------
# TURNS sheet - season cumulative turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 7
$turnsWs.Range("C2").Value = 5
$turnsWs.Range("E2").Value = 13
$turnsWs.Range("D3").Value = 8

# ---------------------------------------------------------------------
# PEN sheet - season cumulative penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 18
$penWs.Range("B3").Value = 18
$penWs.Range("D4").Value = 19
